# chore: update Sheets via scheduled runner
# Refreshes cached market-board figures (currentAveragePrice*, Leve* price/profit
# columns H:N) for a handful of leve rows across several job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 107
$ws.Range("H107").Value = 1429
$ws.Range("I107").Value = 2091.6667
$ws.Range("J107").Value = 932
$ws.Range("K107").Value = 2091.6667
$ws.Range("L107").Value = 932
$ws.Range("M107").Value = -171.6667000000002
$ws.Range("N107").Value = -4772
# Row 114
$ws.Range("H114").Value = 44497
$ws.Range("J114").Value = 44497
$ws.Range("L114").Value = 44497
$ws.Range("N114").Value = -53175
# Row 127
$ws.Range("H127").Value = 1151
$ws.Range("I127").Value = 416.22223
$ws.Range("J127").Value = 2804.25
$ws.Range("K127").Value = 1248.66669
$ws.Range("L127").Value = 8412.75
$ws.Range("M127").Value = 3711.33331
$ws.Range("N127").Value = -18332.75
# Row 129
$ws.Range("H129").Value = 801.1429000000001
$ws.Range("I129").Value = 378.8
$ws.Range("J129").Value = 933.125
$ws.Range("K129").Value = 1136.4
$ws.Range("L129").Value = 2799.375
$ws.Range("M129").Value = 3863.6
$ws.Range("N129").Value = -12799.375
# Row 137
$ws.Range("H137").Value = 18519874
$ws.Range("I137").Value = 924.29266
$ws.Range("J137").Value = 76925790
$ws.Range("K137").Value = 2772.87798
$ws.Range("L137").Value = 230777370
$ws.Range("M137").Value = -222.8779799999998
$ws.Range("N137").Value = -230782470
# Row 141
$ws.Range("H141").Value = 2665
$ws.Range("I141").Value = 1022.5
$ws.Range("K141").Value = 3067.5
$ws.Range("M141").Value = 2112.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8504.947
$ws.Range("I32").Value = 6564.6
$ws.Range("J32").Value = 19592.643
$ws.Range("K32").Value = 6564.6
$ws.Range("L32").Value = 19592.643
$ws.Range("M32").Value = -6277.6
$ws.Range("N32").Value = -20166.643
# Row 45
$ws.Range("H45").Value = 1522.7693
$ws.Range("I45").Value = 1289.5
$ws.Range("J45").Value = 1896
$ws.Range("K45").Value = 1289.5
$ws.Range("L45").Value = 1896
$ws.Range("M45").Value = -912.5
$ws.Range("N45").Value = -2650
# Row 61
$ws.Range("H61").Value = 3834621.8
$ws.Range("I61").Value = 5294688
$ws.Range("K61").Value = 5294688
$ws.Range("M61").Value = -5294476
# Row 74
$ws.Range("H74").Value = 19237126
$ws.Range("I74").Value = 27778664
$ws.Range("J74").Value = 18667.875
$ws.Range("K74").Value = 27778664
$ws.Range("L74").Value = 18667.875
$ws.Range("M74").Value = -27777790
$ws.Range("N74").Value = -20415.875
# Row 77
$ws.Range("H77").Value = 19237126
$ws.Range("I77").Value = 27778664
$ws.Range("J77").Value = 18667.875
$ws.Range("K77").Value = 138893320
$ws.Range("L77").Value = 93339.375
$ws.Range("M77").Value = -138888952
$ws.Range("N77").Value = -102075.375
# Row 122
$ws.Range("H122").Value = 2043.1111
$ws.Range("I122").Value = 1906.2
$ws.Range("J122").Value = 2214.25
$ws.Range("K122").Value = 5718.6
$ws.Range("L122").Value = 6642.75
$ws.Range("M122").Value = -3268.6
$ws.Range("N122").Value = -11542.75
# Row 132
$ws.Range("H132").Value = 774216.3
$ws.Range("I132").Value = 1017002
$ws.Range("J132").Value = 86323.5
$ws.Range("K132").Value = 3051006
$ws.Range("L132").Value = 258970.5
$ws.Range("M132").Value = -3048476
$ws.Range("N132").Value = -264030.5
# Row 136
$ws.Range("H136").Value = 3834621.8
$ws.Range("I136").Value = 5294688
$ws.Range("K136").Value = 15884064
$ws.Range("M136").Value = -15881514

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1374.2368
$ws.Range("I31").Value = 976.3939
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 976.3939
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = -681.3939
$ws.Range("N31").Value = -4590
# Row 34
$ws.Range("H34").Value = 1374.2368
$ws.Range("I34").Value = 976.3939
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 976.3939
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = -774.3939
$ws.Range("N34").Value = -4404

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 20000554
$ws.Range("I5").Value = 515.75
$ws.Range("J5").Value = 55556180
$ws.Range("K5").Value = 1547.25
$ws.Range("L5").Value = 166668540
$ws.Range("M5").Value = -1435.25
$ws.Range("N5").Value = -166668764
# Row 64
$ws.Range("H64").Value = 1078387.2
$ws.Range("I64").Value = 1162
$ws.Range("K64").Value = 3486
$ws.Range("M64").Value = -3216
# Row 67
$ws.Range("H67").Value = 1078387.2
$ws.Range("I67").Value = 1162
$ws.Range("K67").Value = 3486
$ws.Range("M67").Value = -2550
# Row 70
$ws.Range("H70").Value = 28227
$ws.Range("I70").Value = 36857.8
$ws.Range("K70").Value = 110573.4
$ws.Range("M70").Value = -110258.4
# Row 73
$ws.Range("H73").Value = 28227
$ws.Range("I73").Value = 36857.8
$ws.Range("K73").Value = 110573.4
$ws.Range("M73").Value = -109481.4
# Row 87
$ws.Range("H87").Value = 1863.3334
$ws.Range("I87").Value = 1863.3334
$ws.Range("K87").Value = 5590.0002
$ws.Range("M87").Value = -4342.0002
# Row 90
$ws.Range("H90").Value = 1863.3334
$ws.Range("I90").Value = 1863.3334
$ws.Range("K90").Value = 16770.0006
$ws.Range("M90").Value = -10530.0006
# Row 110
$ws.Range("H110").Value = 3130.8572
$ws.Range("I110").Value = 2683.2
$ws.Range("K110").Value = 8049.599999999999
$ws.Range("M110").Value = -3959.599999999999
# Row 131
$ws.Range("H131").Value = 5241.1304
$ws.Range("I131").Value = 5458.9
$ws.Range("J131").Value = 5180.6387
$ws.Range("K131").Value = 16376.7
$ws.Range("L131").Value = 15541.9161
$ws.Range("M131").Value = -11336.7
$ws.Range("N131").Value = -25621.9161
# Row 132
$ws.Range("H132").Value = 47619984
$ws.Range("I132").Value = 200000670
$ws.Range("J132").Value = 1017.125
$ws.Range("K132").Value = 1800006030
$ws.Range("L132").Value = 9154.125
$ws.Range("M132").Value = -1800003500
$ws.Range("N132").Value = -14214.125
# Row 135
$ws.Range("H135").Value = 20000554
$ws.Range("I135").Value = 515.75
$ws.Range("J135").Value = 55556180
$ws.Range("K135").Value = 4641.75
$ws.Range("L135").Value = 500005620
$ws.Range("M135").Value = -2106.75
$ws.Range("N135").Value = -500010690
# Row 137
$ws.Range("H137").Value = 24270.354
$ws.Range("I137").Value = 3882.5
$ws.Range("J137").Value = 28063.441
$ws.Range("K137").Value = 11647.5
$ws.Range("L137").Value = 84190.323
$ws.Range("M137").Value = -6547.5
$ws.Range("N137").Value = -94390.323

$ws = $wb.Worksheets.Item("GSM")
# Row 32
$ws.Range("H32").Value = 24347.5
$ws.Range("J32").Value = 24347.5
$ws.Range("L32").Value = 24347.5
$ws.Range("N32").Value = -24939.5
# Row 80
$ws.Range("H80").Value = 3064.2856
$ws.Range("I80").Value = 3064.2856
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3064.2856
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2066.2856
$ws.Range("N80").ClearContents()
# Row 83
$ws.Range("H83").Value = 3064.2856
$ws.Range("I83").Value = 3064.2856
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 15321.428
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -10329.428
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 3613.5
$ws.Range("I136").Value = 1981.6
$ws.Range("J136").Value = 6333.3335
$ws.Range("K136").Value = 5944.799999999999
$ws.Range("L136").Value = 19000.0005
$ws.Range("M136").Value = -3394.799999999999
$ws.Range("N136").Value = -24100.0005
